# Rewrite the counters_summary sheet: new header layout (B:U) replacing the
# old (B:P) layout, and reset all data-row counter values to 0 to match the
# freshly generated metadata output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (row 1), columns B through U.
$headers = @(
    "COMPLETENESSMANDATORY",
    "COMPLETENESSMANDATORY SCORE",
    "COMPLETENESSOPTIONAL",
    "COMPLETENESSOPTIONAL SCORE",
    "PRECISION",
    "PRECISION SCORE",
    "BUSINESSRULECOMPLIANCE",
    "BUSINESSRULECOMPLIANCE SCORE",
    "METADATACOMPLIANCE",
    "METADATACOMPLIANCE SCORE",
    "UNIQUENESS",
    "UNIQUENESS SCORE",
    "NONREDUNDANCY",
    "NONREDUNDANCY SCORE",
    "SEMANTICCONSISTENCY",
    "SEMANTICCONSISTENCY SCORE",
    "VALUECONSISTENCY",
    "VALUECONSISTENCY SCORE",
    "FORMATCONSISTENCY",
    "FORMATCONSISTENCY SCORE"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    # Column 2 = B, ... Column 21 = U
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Number of data rows currently present below the header (row 2..lastRow).
$lastRow = $ws.UsedRange.Rows.Count

# Reset every counter cell (columns B..U) for every data row to 0, matching
# the regenerated (currently all-zero) metadata counters.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = 0
    }
}
